$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff: Price (D) and Volume(1h) (E) columns
# refreshed for the 22-1-2023 20:xx snapshot. Values are kept as text (matching
# the original inlineStr cells) by forcing a text NumberFormat before assignment
# and then restoring the default "Normal" style so no stray formatting is left behind.
$data = @{
    2 = @{ D="304.20"; E="0.06%" }
    3 = @{ D="36.79"; E="3.35%" }
    4 = @{ D="5.016"; E="-1.58%" }
    5 = @{ D="0.07779"; E="-0.44%" }
    6 = @{ D="2.120"; E="-6.47%" }
    7 = @{ D="8.022"; E="-1.29%" }
    8 = @{ D="0.9204"; E="-0.82%" }
    9 = @{ D="0.09910"; E="2.60%" }
    10 = @{ D="0.1863"; E="2.06%" }
    11 = @{ D="0.08675"; E="-0.58%" }
    12 = @{ D="0.03585"; E="4.78%" }
    13 = @{ D="0.09986"; E="0.34%" }
    14 = @{ D="0.001492"; E="0.13%" }
    15 = @{ D="0.005693"; E="-0.70%" }
    16 = @{ D="3.461"; E="-0.48%" }
    17 = @{ D="4.053"; E="0.75%" }
    18 = @{ D="2.428"; E="13.09%" }
    19 = @{ E="-0.38%" }
    20 = @{ D="0.1305"; E="-1.29%" }
    21 = @{ D="4.939"; E="7.70%" }
    22 = @{ E="-0.99%" }
    23 = @{ D="0.04613"; E="-1.48%" }
    24 = @{ D="0.005143"; E="13.25%" }
    25 = @{ D="0.001238"; E="-0.29%" }
    26 = @{ D="0.0001407"; E="8.05%" }
    39 = @{ D="0.01793"; E="2.01%" }
    40 = @{ D="0.04682"; E="-0.62%" }
    41 = @{ D="0.007706"; E="-2.93%" }
    42 = @{ D="0.1399"; E="-1.52%" }
    43 = @{ D="0.007633"; E="-4.79%" }
    44 = @{ D="0.002141"; E="-6.65%" }
    45 = @{ D="0.01039"; E="14.08%" }
    46 = @{ D="0.00006334"; E="1.58%" }
    47 = @{ D="0.00000000754"; E="0.44%" }
    48 = @{ D="0.0005829"; E="0.49%" }
    49 = @{ D="33.43"; E="485.32%" }
    50 = @{ D="0.002010"; E="-25.34%" }
    51 = @{ D="0.00002111"; E="0.44%" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col]
        $cell.Style = "Normal"
    }
}
